# Applies the "Commercial Drift Test" + "Sheet1" protocol sheets to the
# Micropipette Presentation workbook, matching the target commit:
#   "added and started protocol for testing"

$wb = $excel.ActiveWorkbook

$driftTest = $wb.Worksheets.Item("Drift Test")

# --- Create the two new sheets in an order that reproduces the target
# sheetId allocation (Sheet1 -> sheetId 7, Commercial Drift Test -> sheetId 8)
# while still ending up positioned Drift Test, Commercial Drift Test, Sheet1, ...
$sheet1New = $wb.Worksheets.Add($null, $driftTest)

$cdt = $wb.Worksheets.Add($null, $driftTest)
$cdt.Name = "Commercial Drift Test"

# --- Populate "Commercial Drift Test" ---
$cdt = $wb.Worksheets.Item("Commercial Drift Test")

$cdt.Cells.Item(1, 1).Value = "Trial"
$cdt.Cells.Item(1, 3).Value = "Amount Dispensed Experimental (1mL)Pipette (g)"
$cdt.Cells.Item(1, 4).Value = "Normal Distribution"

$cVals = @(0.2024,0.20469999999999999,0.2039,0.2046,0.20369999999999999,0.2044,0.20469999999999999,0.2044,0.20300000000000001,0.2031,0.20399999999999999,0.2041,0.2049,0.2041,0.2039,0.20419999999999999,0.20430000000000001,0.20380000000000001,0.20430000000000001,0.20380000000000001,0.2039,0.20380000000000001,0.2036,0.2041,0.20319999999999999,0.20269999999999999,0.2039,0.20399999999999999,0.2031,0.20419999999999999,0.2039,0.20380000000000001,0.2036,0.2046,0.2034,0.20319999999999999,0.20369999999999999,0.2034,0.2034,0.20419999999999999,0.2041,0.2039,0.20380000000000001,0.2034,0.20449999999999999,0.20330000000000001,0.2039,0.20369999999999999,0.2041,0.2039)

for ($i = 0; $i -lt $cVals.Length; $i++) {
    $row = $i + 2
    $cdt.Cells.Item($row, 1).Value = $i + 1
    $cdt.Cells.Item($row, 3).Value = $cVals[$i]
    $cdt.Cells.Item($row, 4).Formula = "=NORMDIST(C:C,C54,C55,TRUE )"
}

$cdt.Cells.Item(54, 1).Value = "Mean "
$cdt.Cells.Item(54, 3).Formula = "=AVERAGE(C2:C51)"

$cdt.Cells.Item(55, 1).Value = "Standard Deviation"
$cdt.Cells.Item(55, 3).Formula = "=STDEV(C2:C51)"

$cdt.Cells.Item(56, 1).Value = "Standard Error"
$cdt.Cells.Item(56, 3).Formula = "=(C55/(SQRT(50)))"

$cdt.Columns.Item(3).ColumnWidth = 39
$cdt.Columns.Item(4).ColumnWidth = 28.86

$cdt.Activate()
$cdt.Range("G7").Select()
